# "add headers for sorting"
# The ip_address column on the sort_test.go sheet becomes an "amount" column
# of numeric values (with the last row entered/formatted as a Euro currency
# value), and the selections on both sheets are updated.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("sheet_test.go")
$ws2 = $wb.Worksheets.Item("sort_test.go")

# --- rename the header and replace the ip_address values with amounts ---
$ws2.Cells.Item(1, 6).Value = "amount"

$amounts = @(9, 80, 700, 6000, 50000, 400000, 3000000, 20000000, 100000000)
for ($i = 0; $i -lt $amounts.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 6).Value = $amounts[$i]
}

# last data row (F11) is entered as a literal Euro-formatted value
$ws2.Cells.Item(11, 6).NumberFormat = "@"
$ws2.Cells.Item(11, 6).Value = "€ 9"
$ws2.Cells.Item(11, 6).NumberFormat = "[$€-413]\ #,##0.00;[RED][$€-413]\ #,##0.00\-"

# --- restore/update the selections on both sheets, sort_test.go active ---
[void]$ws2.Activate()
[void]$ws2.Range("F12").Select()
[void]$ws1.Range("A2").Select()
[void]$ws2.Activate()
